$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alunos")
$ws.Activate()
